# Edit: refresh 'Last Updated' timestamp and reshuffle Top Gainers / Top Losers rankings
# per the 2025-10-29 18:25 data pull.

$wb = $excel.ActiveWorkbook

# --- Metadata: bump 'Last Updated' timestamp ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("A2").Value = "29 Oct 2025, 06:25 PM"

# --- Top Gainers: rows 61-76 refreshed with new ranking data ---
$gainers = $wb.Worksheets.Item("Top Gainers")

$gainers.Cells.Item(61, 2).Value = "GREENLAM"
$gainers.Cells.Item(61, 3).Value = 3.8946
$gainers.Cells.Item(61, 4).Value = 3.5858
$gainers.Cells.Item(61, 5).Value = 10.721

$gainers.Cells.Item(62, 2).Value = "APARINDS"
$gainers.Cells.Item(62, 3).Value = 3.8924
$gainers.Cells.Item(62, 4).Value = 8.3414
$gainers.Cells.Item(62, 5).Value = 15.5876

$gainers.Cells.Item(63, 2).Value = "HITECHGEAR"
$gainers.Cells.Item(63, 3).Value = 3.8587
$gainers.Cells.Item(63, 4).Value = 1.1486
$gainers.Cells.Item(63, 5).Value = 9.9254

$gainers.Cells.Item(64, 2).Value = "ORIENTTECH"
$gainers.Cells.Item(64, 3).Value = 3.827
$gainers.Cells.Item(64, 4).Value = 0.5247
$gainers.Cells.Item(64, 5).Value = 32.6784

$gainers.Cells.Item(65, 2).Value = "ICRA"
$gainers.Cells.Item(65, 3).Value = 3.7985
$gainers.Cells.Item(65, 4).Value = 4.4793
$gainers.Cells.Item(65, 5).Value = 2.8828

$gainers.Cells.Item(66, 2).Value = "SALASAR"
$gainers.Cells.Item(66, 3).Value = 3.7935
$gainers.Cells.Item(66, 4).Value = 4.7872
$gainers.Cells.Item(66, 5).Value = 11.0485

$gainers.Cells.Item(67, 2).Value = "NPST"
$gainers.Cells.Item(67, 3).Value = 3.7841
$gainers.Cells.Item(67, 4).Value = -2.0689
$gainers.Cells.Item(67, 5).Value = -3.5677

$gainers.Cells.Item(68, 2).Value = "DCW"
$gainers.Cells.Item(68, 3).Value = 3.7544
$gainers.Cells.Item(68, 4).Value = 2.3219
$gainers.Cells.Item(68, 5).Value = -3.9753

$gainers.Cells.Item(69, 2).Value = "RHETAN"
$gainers.Cells.Item(69, 3).Value = 3.754
$gainers.Cells.Item(69, 4).Value = 4.178
$gainers.Cells.Item(69, 5).Value = 6.549

$gainers.Cells.Item(70, 2).Value = "HINDPETRO"
$gainers.Cells.Item(70, 3).Value = 3.6935
$gainers.Cells.Item(70, 4).Value = 6.9335
$gainers.Cells.Item(70, 5).Value = 5.7397

$gainers.Cells.Item(71, 2).Value = "BHARTIHEXA"
$gainers.Cells.Item(71, 3).Value = 3.6718
$gainers.Cells.Item(71, 4).Value = 7.0877
$gainers.Cells.Item(71, 5).Value = 15.3332

$gainers.Cells.Item(72, 2).Value = "HLEGLAS"
$gainers.Cells.Item(72, 3).Value = 3.659
$gainers.Cells.Item(72, 4).Value = 8.1155
$gainers.Cells.Item(72, 5).Value = 27.1239

$gainers.Cells.Item(73, 2).Value = "RHIM"
$gainers.Cells.Item(73, 3).Value = 3.6544
$gainers.Cells.Item(73, 4).Value = 3.2276
$gainers.Cells.Item(73, 5).Value = 5.1826

$gainers.Cells.Item(74, 2).Value = "SHK"
$gainers.Cells.Item(74, 3).Value = 3.6347
$gainers.Cells.Item(74, 4).Value = 2.388
$gainers.Cells.Item(74, 5).Value = -1.932

$gainers.Cells.Item(75, 2).Value = "BCLIND"
$gainers.Cells.Item(75, 3).Value = 3.6271
$gainers.Cells.Item(75, 4).Value = 2.2945
$gainers.Cells.Item(75, 5).Value = 0.1728

$gainers.Cells.Item(76, 2).Value = "MUKANDLTD"
$gainers.Cells.Item(76, 3).Value = 3.6133
$gainers.Cells.Item(76, 4).Value = 11.9685
$gainers.Cells.Item(76, 5).Value = 9.5508

# --- Top Losers: two isolated weekly-% corrections ---
$losers = $wb.Worksheets.Item("Top Losers")
$losers.Cells.Item(18, 4).Value = -0.062   # CRAMC
$losers.Cells.Item(48, 4).Value = 0.05     # RUBICON

# --- Top Losers: rows 51-73 refreshed with new ranking data ---
$losers.Cells.Item(51, 2).Value = "UNIMECH"
$losers.Cells.Item(51, 3).Value = -2.8008
$losers.Cells.Item(51, 4).Value = -1.6104
$losers.Cells.Item(51, 5).Value = -0.4585

$losers.Cells.Item(52, 2).Value = "TTKPRESTIG"
$losers.Cells.Item(52, 3).Value = -2.7438
$losers.Cells.Item(52, 4).Value = 8.0012
$losers.Cells.Item(52, 5).Value = 9.6505

$losers.Cells.Item(53, 2).Value = "PFOCUS"
$losers.Cells.Item(53, 3).Value = -2.7039
$losers.Cells.Item(53, 4).Value = -2.6276
$losers.Cells.Item(53, 5).Value = -1.2163

$losers.Cells.Item(54, 2).Value = "ALLDIGI"
$losers.Cells.Item(54, 3).Value = -2.6342
$losers.Cells.Item(54, 4).Value = -0.2306
$losers.Cells.Item(54, 5).Value = -5.3103

$losers.Cells.Item(55, 2).Value = "PRIVISCL"
$losers.Cells.Item(55, 3).Value = -2.6288
$losers.Cells.Item(55, 4).Value = -2.1048
$losers.Cells.Item(55, 5).Value = 19.7451

$losers.Cells.Item(56, 2).Value = "CANHLIFE"
$losers.Cells.Item(56, 3).Value = -2.6148
$losers.Cells.Item(56, 4).Value = 3.7771
$losers.Cells.Item(56, 5).Value = "N/A"

$losers.Cells.Item(57, 2).Value = "GKENERGY"
$losers.Cells.Item(57, 3).Value = -2.6122
$losers.Cells.Item(57, 4).Value = -9.8077
$losers.Cells.Item(57, 5).Value = 23.2758

$losers.Cells.Item(58, 2).Value = "SGFIN"
$losers.Cells.Item(58, 3).Value = -2.592
$losers.Cells.Item(58, 4).Value = -0.0627
$losers.Cells.Item(58, 5).Value = 11.7235

$losers.Cells.Item(59, 2).Value = "ARVINDFASN"
$losers.Cells.Item(59, 3).Value = -2.549
$losers.Cells.Item(59, 4).Value = -2.9892
$losers.Cells.Item(59, 5).Value = -4.4223

$losers.Cells.Item(60, 2).Value = "EDELWEISS"
$losers.Cells.Item(60, 3).Value = -2.5422
$losers.Cells.Item(60, 4).Value = -3.3745
$losers.Cells.Item(60, 5).Value = 8.5305

$losers.Cells.Item(61, 2).Value = "SAMHI"
$losers.Cells.Item(61, 3).Value = -2.5284
$losers.Cells.Item(61, 4).Value = 1.8231
$losers.Cells.Item(61, 5).Value = 2.8516

$losers.Cells.Item(62, 2).Value = "TBOTEK"
$losers.Cells.Item(62, 3).Value = -2.524
$losers.Cells.Item(62, 4).Value = -3.5732
$losers.Cells.Item(62, 5).Value = 1.036

$losers.Cells.Item(63, 2).Value = "UJJIVANSFB"
$losers.Cells.Item(63, 3).Value = -2.5201
$losers.Cells.Item(63, 4).Value = 0.3845
$losers.Cells.Item(63, 5).Value = 12.6645

$losers.Cells.Item(64, 2).Value = "AMBER"
$losers.Cells.Item(64, 3).Value = -2.5098
$losers.Cells.Item(64, 4).Value = -0.1082
$losers.Cells.Item(64, 5).Value = 2.763

$losers.Cells.Item(65, 2).Value = "GRPLTD"
$losers.Cells.Item(65, 3).Value = -2.4898
$losers.Cells.Item(65, 4).Value = -5.9894
$losers.Cells.Item(65, 5).Value = -5.4586

$losers.Cells.Item(66, 2).Value = "NESCO"
$losers.Cells.Item(66, 3).Value = -2.4722
$losers.Cells.Item(66, 4).Value = 1.9934
$losers.Cells.Item(66, 5).Value = 3.8931

$losers.Cells.Item(67, 2).Value = "PILANIINVS"
$losers.Cells.Item(67, 3).Value = -2.4546
$losers.Cells.Item(67, 4).Value = -0.7907
$losers.Cells.Item(67, 5).Value = 4.267

$losers.Cells.Item(68, 2).Value = "NSIL"
$losers.Cells.Item(68, 3).Value = -2.4088
$losers.Cells.Item(68, 4).Value = -1.7646
$losers.Cells.Item(68, 5).Value = 4.7431

$losers.Cells.Item(69, 2).Value = "COALINDIA"
$losers.Cells.Item(69, 3).Value = -2.4016
$losers.Cells.Item(69, 4).Value = -3.058
$losers.Cells.Item(69, 5).Value = -2.0387

$losers.Cells.Item(70, 2).Value = "JNKINDIA"
$losers.Cells.Item(70, 3).Value = -2.3482
$losers.Cells.Item(70, 4).Value = -2.8371
$losers.Cells.Item(70, 5).Value = 4.2622

$losers.Cells.Item(72, 2).Value = "DEEDEV"
$losers.Cells.Item(72, 3).Value = -2.3334
$losers.Cells.Item(72, 4).Value = -6.6528
$losers.Cells.Item(72, 5).Value = -7.4227

$losers.Cells.Item(73, 2).Value = "WEALTH"
$losers.Cells.Item(73, 3).Value = -2.2793
$losers.Cells.Item(73, 4).Value = -2.7981
$losers.Cells.Item(73, 5).Value = -2.7981

Write-Host "Applied 2025-10-29 18:25 data refresh."